$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OutputForces")

# Updated buckling safety factor: the buckle load formula now divides the
# effective modulus by the buckle safety factor (C20) as well.
$ws.Range("G21:L21").Formula = "=-G14*(G19/1000)^2/(PI()^2*`$C`$17/`$C`$20)"

# Selection left where the user last clicked while making this edit.
[void]$ws.Range("F33").Select()
